$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new date
$ws.Name = "Through 2022-06-24"

# Update the June label text to reflect the new date
$ws.Range("A7").Value = "June (through 06-24)"

# Update June row (row 7) values
$ws.Range("B7").Value = 14
$ws.Range("C7").Value = 29
$ws.Range("D7").Value = 59
$ws.Range("E7").Value = 46
$ws.Range("F7").Value = 35
$ws.Range("G7").Value = 92
$ws.Range("H7").Value = 94
$ws.Range("I7").Value = 114

# Update Total row (row 8) values
$ws.Range("B8").Value = 122
$ws.Range("C8").Value = 238
$ws.Range("D8").Value = 375
$ws.Range("E8").Value = 341
$ws.Range("F8").Value = 239
$ws.Range("G8").Value = 450
$ws.Range("H8").Value = 725
$ws.Range("I8").Value = 777
